# Update the cryptocurrency price/volume data on the active sheet.
# Values are written with a leading apostrophe to force Excel to store
# them as text (preserving formats such as "37.384.68" or "4.40"
# instead of auto-converting them to numbers), and each cell's style
# is restored afterwards so no incidental formatting changes are made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "37.384.68"
    "E2" = "  +0.11%  "
    "D3" = "2.065.47"
    "E3" = "  -0.51%  "
    "E4" = "  +0.05%  "
    "D5" = "234.13"
    "E5" = "  -0.21%  "
    "E6" = "  +1.58%  "
    "E7" = "  -0.02%  "
    "D8" = "56.97"
    "E8" = "  -0.64%  "
    "D10" = "0.0762"
    "E10" = "  +0.37%  "
    "E11" = "  +0.65%  "
    "D12" = "2.370.13"
    "E12" = "  -0.46%  "
    "D13" = "14.61"
    "E13" = "  +0.79%  "
    "E14" = "  -1.75%  "
    "D15" = "0.777"
    "E15" = "  +0.44%  "
    "D16" = "5.13"
    "E16" = "  -2.00%  "
    "D17" = "2.064.38"
    "E17" = "  -0.22%  "
    "D18" = "37.305.98"
    "E18" = "  -0.45%  "
    "E19" = "  +4.69%  "
    "E20" = "  +1.60%  "
    "D21" = "0.0₃0810"
    "E21" = "  -0.05%  "
    "D22" = "226.45"
    "E22" = "  +1.31%  "
    "E23" = "  +0.01%  "
    "D24" = "2.45"
    "E24" = "  +0.97%  "
    "D25" = "2.40"
    "E25" = "  -0.96%  "
    "D26" = "166.56"
    "E26" = "  +2.38%  "
    "D27" = "8.79"
    "E27" = "  -0.91%  "
    "E28" = "  +4.24%  "
    "D29" = "19.07"
    "E29" = "  -0.96%  "
    "E30" = "  -3.15%  "
    "E31" = "  +0.05%  "
    "D32" = "4.48"
    "E32" = "  +0.79%  "
    "E33" = "  -1.06%  "
    "E34" = "  +3.81%  "
    "E35" = "  -2.61%  "
    "E36" = "  +0.01%  "
    "E37" = "  +0.05%  "
    "E38" = "  -2.20%  "
    "D39" = "5.67"
    "E39" = "  -4.95%  "
    "D40" = "2.96"
    "E40" = "  -0.23%  "
    "B41" = "FTXToken"
    "C41" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D41" = "4.40"
    "E41" = "  +1.05%  "
    "D42" = "1.466.04"
    "E42" = "  -0.44%  "
    "B43" = "Aave"
    "C43" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D43" = "96.18"
    "E43" = "  +1.15%  "
    "D44" = "0.0939"
    "E44" = "  -2.15%  "
    "E45" = "  +3.79%  "
    "E46" = "  +1.38%  "
    "E47" = "  -1.14%  "
    "D48" = "15.03"
    "E48" = "  -6.91%  "
    "D49" = "7.14"
    "E49" = "  -1.89%  "
    "E50" = "  +0.65%  "
    "D51" = "2.258.59"
    "E51" = "  -0.43%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $savedStyle = $cell.Style
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = $savedStyle
}

Write-Host ("Updated {0} cells" -f $updates.Count)
